$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the paragraph that currently reads "ddddd" (the second paragraph).
$target = $d.Paragraphs.Item(2)
$r = $target.Range

# Replace its content with "D" + "dddd" split across two runs, wrapped in a
# spell-check proofErr pair (capitalization "correction"), drop its old
# bookmark, add a new blank paragraph right after it, and append a further
# paragraph - underlined - holding the "A numaaaaaaaaaaa" comment, carrying
# forward the _GoBack bookmark.
$xml = "<w:p $W>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:t>D</w:t></w:r>" + `
    "<w:r><w:t>dddd</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>" + `
  "</w:p>" + `
  "<w:p $W/>" + `
  "<w:p $W>" + `
    "<w:pPr><w:rPr><w:u w:val='single'/></w:rPr></w:pPr>" + `
    "<w:r><w:t>A numaaaaaaaaaaa</w:t></w:r>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" + `
    "<w:bookmarkEnd w:id='0'/>" + `
  "</w:p>"

$r.InsertXML($xml)
